# Updated amazon ubuntu14 ami's
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = "ami-4a9b1930"
$ws.Range("G4").Value = "ami-1d90a97d"
$ws.Range("H4").Value = "ami-0fcf1c77"
$ws.Range("I4").Value = "ami-308d2749"
$ws.Range("J4").Value = "ami-5029a93f"
$ws.Range("K4").Value = "ami-3ec8d75a"

$ws.Range("H9").Select()
